# Update Summary Info report with refreshed source numbers.
# Only the raw input columns (Current Month Active, Prior Month Debit,
# Current Active DVH, Lifetime Cancels DVH, Current Active Copay) change;
# the formula columns (E,H,I,K,O,P,T,U) recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 320
$ws.Range("M2").Value = 36
$ws.Range("R2").Value = 9

$ws.Range("C3").Value = 976
$ws.Range("J3").Value = 380
$ws.Range("M3").Value = 109

$ws.Range("C4").Value = 636
$ws.Range("J4").Value = 22
$ws.Range("M4").Value = 127
$ws.Range("R4").Value = 11

$ws.Range("C5").Value = 301
$ws.Range("M5").Value = 52
$ws.Range("R5").Value = 16

$ws.Range("C6").Value = 142
$ws.Range("J6").Value = 1
$ws.Range("M6").Value = 26

$ws.Range("C7").Value = 958
$ws.Range("J7").Value = 377
$ws.Range("M7").Value = 245
$ws.Range("N7").Value = 2
$ws.Range("R7").Value = 17

$ws.Range("C8").Value = 928
$ws.Range("J8").Value = 219
$ws.Range("M8").Value = 129

$ws.Range("C9").Value = 386
$ws.Range("M9").Value = 90

$ws.Range("C10").Value = 489
$ws.Range("J10").Value = 28
$ws.Range("M10").Value = 85
$ws.Range("R10").Value = 2

$ws.Range("C11").Value = 450
$ws.Range("J11").Value = 0
$ws.Range("M11").Value = 105
$ws.Range("R11").Value = 22

$ws.Range("C12").Value = 398
$ws.Range("M12").Value = 73
$ws.Range("R12").Value = 5

$ws.Range("C13").Value = 110
$ws.Range("M13").Value = 24
$ws.Range("R13").Value = 1

$ws.Range("C14").Value = 150
$ws.Range("M14").Value = 16

$ws.Range("C15").Value = 734
$ws.Range("J15").Value = 90
$ws.Range("M15").Value = 132

$ws.Range("C16").Value = 881
$ws.Range("J16").Value = 394
$ws.Range("M16").Value = 116
$ws.Range("R16").Value = 24

$ws.Range("C17").Value = 586
$ws.Range("J17").Value = 41
$ws.Range("M17").Value = 120

$ws.Range("C18").Value = 748
$ws.Range("J18").Value = 146
$ws.Range("M18").Value = 141

$ws.Range("C19").Value = 622
$ws.Range("J19").Value = 56
$ws.Range("M19").Value = 107

$ws.Range("C20").Value = 505
$ws.Range("J20").Value = 32
$ws.Range("M20").Value = 80

$ws.Range("C21").Value = 1143
$ws.Range("J21").Value = 518
$ws.Range("M21").Value = 98
$ws.Range("R21").Value = 3

$ws.Range("C22").Value = 635
$ws.Range("J22").Value = 322
$ws.Range("M22").Value = 62

$ws.Range("J23").Value = 124

$ws.Range("C24").Value = 388
$ws.Range("M24").Value = 82
$ws.Range("R24").Value = 20
